$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 55557840
$ws.Range("I40").Value = 1665.3334
$ws.Range("J40").Value = 83335930
$ws.Range("K40").Value = 1665.3334
$ws.Range("L40").Value = 83335930
$ws.Range("M40").Value = -1490.3334
$ws.Range("N40").Value = -83336280
$ws.Range("H43").Value = 5426.7856
$ws.Range("I43").Value = 5398.4
$ws.Range("J43").Value = 5497.75
$ws.Range("K43").Value = 5398.4
$ws.Range("L43").Value = 5497.75
$ws.Range("M43").Value = -5329.4
$ws.Range("N43").Value = -5635.75
$ws.Range("H51").Value = 10985.9
$ws.Range("J51").Value = 6617.6665
$ws.Range("L51").Value = 6617.6665
$ws.Range("N51").Value = -7585.6665
$ws.Range("H98").Value = 1873.45
$ws.Range("I98").Value = 1270.25
$ws.Range("J98").Value = 4286.25
$ws.Range("K98").Value = 1270.25
$ws.Range("L98").Value = 4286.25
$ws.Range("M98").Value = 227.75
$ws.Range("N98").Value = -7282.25
$ws.Range("H103").Value = 41668770
$ws.Range("J103").Value = 41668770
$ws.Range("L103").Value = 125006310
$ws.Range("N103").Value = -125007482
$ws.Range("H113").Value = 2782.3794
$ws.Range("I113").Value = 2866.6843
$ws.Range("J113").Value = 2622.2
$ws.Range("K113").Value = 2866.6843
$ws.Range("L113").Value = 2622.2
$ws.Range("M113").Value = 387.3157000000001
$ws.Range("N113").Value = -9130.200000000001
$ws.Range("H122").Value = 1873.45
$ws.Range("I122").Value = 1270.25
$ws.Range("J122").Value = 4286.25
$ws.Range("K122").Value = 3810.75
$ws.Range("L122").Value = 12858.75
$ws.Range("M122").Value = -1360.75
$ws.Range("N122").Value = -17758.75
$ws.Range("H125").Value = 987.2778
$ws.Range("I125").Value = 421
$ws.Range("J125").Value = 1347.6364
$ws.Range("K125").Value = 3789
$ws.Range("L125").Value = 12128.7276
$ws.Range("M125").Value = -1329
$ws.Range("N125").Value = -17048.7276
$ws.Range("H127").Value = 2295.6667
$ws.Range("I127").Value = 2295.6667
$ws.Range("K127").Value = 6887.000100000001
$ws.Range("M127").Value = -1927.000100000001
$ws.Range("H129").Value = 6756.875
$ws.Range("I129").Value = 661.8889
$ws.Range("J129").Value = 10413.866
$ws.Range("K129").Value = 1985.6667
$ws.Range("L129").Value = 31241.598
$ws.Range("M129").Value = 3014.3333
$ws.Range("N129").Value = -41241.598
$ws.Range("H131").Value = 4547524
$ws.Range("I131").Value = 583.3333
$ws.Range("K131").Value = 1749.9999
$ws.Range("M131").Value = 3290.0001
$ws.Range("H141").Value = 5385.1113
$ws.Range("I141").Value = 5385.1113
$ws.Range("K141").Value = 16155.3339
$ws.Range("M141").Value = -10975.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2198.158
$ws.Range("I63").Value = 2261.75
$ws.Range("J63").Value = 1859
$ws.Range("K63").Value = 2261.75
$ws.Range("L63").Value = 1859
$ws.Range("M63").Value = -1575.75
$ws.Range("N63").Value = -3231
$ws.Range("H66").Value = 2198.158
$ws.Range("I66").Value = 2261.75
$ws.Range("J66").Value = 1859
$ws.Range("K66").Value = 11308.75
$ws.Range("L66").Value = 9295
$ws.Range("M66").Value = -7876.75
$ws.Range("N66").Value = -16159

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 6380.5
$ws.Range("I11").Value = 1167
$ws.Range("J11").Value = 8987.25
$ws.Range("K11").Value = 1167
$ws.Range("L11").Value = 8987.25
$ws.Range("M11").Value = -1027
$ws.Range("N11").Value = -9267.25
$ws.Range("H20").Value = 5811.7144
$ws.Range("I20").Value = 7207.095
$ws.Range("J20").Value = 3718.6428
$ws.Range("K20").Value = 7207.095
$ws.Range("L20").Value = 3718.6428
$ws.Range("M20").Value = -6960.095
$ws.Range("N20").Value = -4212.6428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2060.7273
$ws.Range("I122").Value = 1458.6
$ws.Range("J122").Value = 2987.077
$ws.Range("K122").Value = 4375.799999999999
$ws.Range("L122").Value = 8961.231
$ws.Range("M122").Value = -1925.799999999999
$ws.Range("N122").Value = -13861.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 38.416668
$ws.Range("J38").Value = 50
$ws.Range("L38").Value = 150
$ws.Range("N38").Value = -844
$ws.Range("H68").Value = 2604.4443
$ws.Range("J68").Value = 4097.4
$ws.Range("L68").Value = 12292.2
$ws.Range("N68").Value = -13914.2
$ws.Range("H71").Value = 2604.4443
$ws.Range("J71").Value = 4097.4
$ws.Range("L71").Value = 36876.6
$ws.Range("N71").Value = -44988.6
$ws.Range("H80").Value = 27787326
$ws.Range("I80").Value = 33344194
$ws.Range("K80").Value = 100032582
$ws.Range("M80").Value = -100031646
$ws.Range("H83").Value = 27787326
$ws.Range("I83").Value = 33344194
$ws.Range("K83").Value = 300097746
$ws.Range("M83").Value = -300093066
$ws.Range("H97").Value = 283.81818
$ws.Range("J97").Value = 247.42857
$ws.Range("L97").Value = 742.28571
$ws.Range("N97").Value = -1734.28571
$ws.Range("H98").Value = 642.75
$ws.Range("I98").Value = 536.75
$ws.Range("K98").Value = 1610.25
$ws.Range("M98").Value = -112.25
$ws.Range("H113").Value = 1552.5454
$ws.Range("I113").Value = 1289.3334
$ws.Range("J113").Value = 1868.4
$ws.Range("K113").Value = 3868.0002
$ws.Range("L113").Value = 5605.200000000001
$ws.Range("M113").Value = -1698.0002
$ws.Range("N113").Value = -9945.200000000001
$ws.Range("H129").Value = 9619012
$ws.Range("I129").Value = 14707551
$ws.Range("K129").Value = 44122653
$ws.Range("M129").Value = -44117653
$ws.Range("H131").Value = 4233.3105
$ws.Range("I131").Value = 3146.1667
$ws.Range("K131").Value = 9438.500100000001
$ws.Range("M131").Value = -4398.500100000001
$ws.Range("H139").Value = 5213.5835
$ws.Range("I139").Value = 3196.76
$ws.Range("K139").Value = 9590.280000000001
$ws.Range("M139").Value = -4450.280000000001
$ws.Range("H140").Value = 3531.394
$ws.Range("I140").Value = 1571.4073
$ws.Range("J140").Value = 12351.333
$ws.Range("K140").Value = 4714.2219
$ws.Range("L140").Value = 37053.999
$ws.Range("M140").Value = 465.7780999999995
$ws.Range("N140").Value = -47413.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 18449.166
$ws.Range("I5").Value = 33398.332
$ws.Range("K5").Value = 33398.332
$ws.Range("M5").Value = -33286.332
$ws.Range("H57").Value = 94999.5
$ws.Range("I57").Value = 90000
$ws.Range("J57").Value = 99999
$ws.Range("K57").Value = 90000
$ws.Range("L57").Value = 99999
$ws.Range("M57").Value = -89180
$ws.Range("N57").Value = -101639
$ws.Range("H70").Value = 4998.4614
$ws.Range("I70").Value = 4980
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 4980
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -4710
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 4998.4614
$ws.Range("I73").Value = 4980
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 4980
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -4044
$ws.Range("N73").Value = -6872
$ws.Range("H107").Value = 362.3
$ws.Range("I107").Value = 204.125
$ws.Range("K107").Value = 204.125
$ws.Range("M107").Value = 1715.875
$ws.Range("H132").Value = 4341961
$ws.Range("I132").Value = 3372.1516
$ws.Range("K132").Value = 10116.4548
$ws.Range("M132").Value = -7586.4548

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2045.1428
$ws.Range("I16").Value = 2912
$ws.Range("K16").Value = 2912
$ws.Range("M16").Value = -2742
$ws.Range("H46").Value = 952.4211
$ws.Range("I46").Value = 808.3333
$ws.Range("K46").Value = 808.3333
$ws.Range("M46").Value = -620.3333
$ws.Range("H68").Value = 1606182.1
$ws.Range("I68").Value = 2454469.8
$ws.Range("J68").Value = 3860.889
$ws.Range("K68").Value = 2454469.8
$ws.Range("L68").Value = 3860.889
$ws.Range("M68").Value = -2453720.8
$ws.Range("N68").Value = -5358.889
$ws.Range("H71").Value = 1606182.1
$ws.Range("I71").Value = 2454469.8
$ws.Range("J71").Value = 3860.889
$ws.Range("K71").Value = 12272349
$ws.Range("L71").Value = 19304.445
$ws.Range("M71").Value = -12268605
$ws.Range("N71").Value = -26792.445
$ws.Range("H82").Value = 2669.9614
$ws.Range("I82").Value = 825.4
$ws.Range("K82").Value = 825.4
$ws.Range("M82").Value = -464.4
$ws.Range("H85").Value = 2669.9614
$ws.Range("I85").Value = 825.4
$ws.Range("K85").Value = 825.4
$ws.Range("M85").Value = 422.6
$ws.Range("H100").Value = 25029694
$ws.Range("I100").Value = 5247.5
$ws.Range("J100").Value = 41712660
$ws.Range("K100").Value = 5247.5
$ws.Range("L100").Value = 41712660
$ws.Range("M100").Value = -4706.5
$ws.Range("N100").Value = -41713742
$ws.Range("H132").Value = 3583.6216
$ws.Range("I132").Value = 2352.375
$ws.Range("K132").Value = 7057.125
$ws.Range("M132").Value = -4527.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 39248.75
$ws.Range("I51").Value = 34535
$ws.Range("J51").Value = 43962.5
$ws.Range("K51").Value = 34535
$ws.Range("L51").Value = 43962.5
$ws.Range("M51").Value = -34025
$ws.Range("N51").Value = -44982.5
$ws.Range("H58").Value = 58333
$ws.Range("J58").Value = 62499.5
$ws.Range("L58").Value = 62499.5
$ws.Range("N58").Value = -63115.5
$ws.Range("H62").Value = 7814.7617
$ws.Range("I62").Value = 4648.8335
$ws.Range("J62").Value = 12036
$ws.Range("K62").Value = 4648.8335
$ws.Range("L62").Value = 12036
$ws.Range("M62").Value = -4024.8335
$ws.Range("N62").Value = -13284
$ws.Range("H65").Value = 7814.7617
$ws.Range("I65").Value = 4648.8335
$ws.Range("J65").Value = 12036
$ws.Range("K65").Value = 23244.1675
$ws.Range("L65").Value = 60180
$ws.Range("M65").Value = -20124.1675
$ws.Range("N65").Value = -66420
$ws.Range("H132").Value = 324277.6
$ws.Range("I132").Value = 1947.9131
$ws.Range("J132").Value = 1250975.4
$ws.Range("K132").Value = 5843.7393
$ws.Range("L132").Value = 3752926.2
$ws.Range("M132").Value = -3313.7393
$ws.Range("N132").Value = -3757986.2
